$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J25").Value = -10.57102434548657
$ws.Range("K25").Value = 3.574792626865744
$ws.Range("I26").Value = -10.51602434548657
$ws.Range("J26").Value = 3.629792626865739
$ws.Range("H27").Value = -10.51602434548659
$ws.Range("I27").Value = 3.629792626865722
$ws.Range("G28").Value = -10.51602434548657
$ws.Range("H28").Value = 3.629792626865739
$ws.Range("F29").Value = -10.51602434548659
$ws.Range("G29").Value = 3.689792626865739
$ws.Range("H29").Value = 5.716794673020033
$ws.Range("I29").Value = -2.932891992481572
$ws.Range("J29").Value = -5.907106464233138
$ws.Range("K29").Value = 2.050717317831413
$ws.Range("E30").Value = -10.56702434548656
$ws.Range("F30").Value = 3.629792626865722
$ws.Range("G30").Value = 5.89679467302004
$ws.Range("H30").Value = -2.792891992481557
$ws.Range("I30").Value = -5.867106464233117
$ws.Range("J30").Value = 2.090717317831433
$ws.Range("D31").Value = -7.516024345486572
$ws.Range("E31").Value = 5.729792626865716
$ws.Range("F31").Value = 7.756794673020053
$ws.Range("G31").Value = -1.632891992481575
$ws.Range("H31").Value = -5.867603096389431
$ws.Range("I31").Value = 2.090222366053069
$ws.Range("C32").Value = -2.61602434548657
$ws.Range("D32").Value = 6.329792626865739
$ws.Range("E32").Value = 6.156794673020033
$ws.Range("F32").Value = -2.232891992481555
$ws.Range("G32").Value = -5.80710646423312
$ws.Range("H32").Value = 2.15071731783143
$ws.Range("B33").Value = -3.56602434548657
$ws.Range("C33").Value = 5.769792626865737
$ws.Range("D33").Value = 6.056794673020022
$ws.Range("E33").Value = -2.332891992481564
$ws.Range("F33").Value = -5.807106464233115
$ws.Range("G33").Value = 2.150717317831436
$ws.Range("H33").Value = 0.5746947653272656
$ws.Range("I33").Value = -0.3780328061469476
$ws.Range("J33").Value = -0.5979284922632784
$ws.Range("K33").Value = 1.443932912058642
$ws.Range("B34").Value = 2.529792626865742
$ws.Range("C34").Value = 5.356794673020033
$ws.Range("D34").Value = -0.1328919924815612
$ws.Range("E34").Value = -5.307106464233129
$ws.Range("F34").Value = 2.39071731783143
$ws.Range("G34").Value = 0.5746947653272656
$ws.Range("H34").Value = -0.3780328061469476
$ws.Range("I34").Value = -0.5979284922632784
$ws.Range("J34").Value = 1.443932912058642
$ws.Range("B35").Value = 2.956794673020027
$ws.Range("C35").Value = 0.06710800751844204
$ws.Range("D35").Value = -1.607106464233111
$ws.Range("E35").Value = 3.550717317831427
$ws.Range("F35").Value = 0.804913431933457
$ws.Range("G35").Value = -0.2518180252671272
$ws.Range("H35").Value = -0.580456746937102
$ws.Range("I35").Value = 1.430045938313753
$ws.Range("B36").Value = 0.2479822279162387
$ws.Range("C36").Value = 0.4928935357668678
$ws.Range("D36").Value = 4.150717317831436
$ws.Range("E36").Value = 0.7346947653272622
$ws.Range("F36").Value = -0.3180328061469595
$ws.Range("G36").Value = -0.5635338037639739
$ws.Range("H36").Value = 1.390442373121999
$ws.Range("B37").Value = 0.1928935357668848
$ws.Range("C37").Value = 4.350717317831439
$ws.Range("D37").Value = 1.034694765327259
$ws.Range("E37").Value = -0.1180328061469567
$ws.Range("F37").Value = -0.4379284922632818
$ws.Range("G37").Value = 1.50393291205863
$ws.Range("H37").Value = 1.429861316022425
$ws.Range("I37").Value = 0.1526454189713746
$ws.Range("J37").Value = 0.4656763841019966
$ws.Range("K37").Value = -0.03673872975578271
$ws.Range("B38").Value = 0.6507173178314358
$ws.Range("C38").Value = -1.265305234672738
$ws.Range("D38").Value = 3.681967193853055
$ws.Range("E38").Value = 1.412071507736727
$ws.Range("F38").Value = 1.523932912058654
$ws.Range("G38").Value = 1.579861316022431
$ws.Range("H38").Value = 0.3026454189713803
$ws.Range("I38").Value = 0.5156763841019938
$ws.Range("J38").Value = 0.0132612702442145
$ws.Range("B39").Value = 0.5556414999948345
$ws.Range("C39").Value = 0.7719671938530577
$ws.Range("D39").Value = 0.1360272089500689
$ws.Range("E39").Value = 2.103932912058638
$ws.Range("F39").Value = 2.199861316022421
$ws.Range("G39").Value = 0.9826454189713727
$ws.Range("H39").Value = 0.9456763841019864
$ws.Range("I39").Value = 0.1832612702442162
$ws.Range("B40").Value = 0.7219671938530607
$ws.Range("C40").Value = 0.162071507736727
$ws.Range("D40").Value = 1.853932912058638
$ws.Range("E40").Value = 1.929861316022411
$ws.Range("F40").Value = 0.5526454189713661
$ws.Range("G40").Value = 0.8656763841019881
$ws.Range("H40").Value = 0.113261270244223
$ws.Range("B41").Value = -1.186928492263277
$ws.Range("C41").Value = 0.3639329120586434
$ws.Range("D41").Value = -0.1101386839775809
$ws.Range("E41").Value = -0.4473545810286339
$ws.Range("F41").Value = 0.5156763841019938
$ws.Range("G41").Value = 0.2632612702442003
$ws.Range("H41").Value = 1.014894517653033
$ws.Range("I41").Value = 0.8219948986110523
$ws.Range("J41").Value = 0.3112534500261574
$ws.Range("K41").Value = 0.5166541826233555
$ws.Range("B42").Value = 0.5039362906915836
$ws.Range("C42").Value = 0.2797613160224302
$ws.Range("D42").Value = -0.3773545810286265
$ws.Range("E42").Value = 0.805676384102
$ws.Range("F42").Value = 0.3032612702442208
$ws.Range("G42").Value = 0.9148945176530248
$ws.Range("H42").Value = 0.5219948986110409
$ws.Range("I42").Value = 0.01125345002614597
$ws.Range("J42").Value = 0.1566541826233561
$ws.Range("B43").Value = 0.4898613160224272
$ws.Range("C43").Value = -0.1073545810286305
$ws.Range("D43").Value = 0.6656763841019995
$ws.Range("E43").Value = 0.3632612702442088
$ws.Range("F43").Value = 0.9148945176530248
$ws.Range("G43").Value = 0.5519948986110705
$ws.Range("H43").Value = -0.04874654997384201
$ws.Range("I43").Value = 0.116654182623364
$ws.Range("B44").Value = -0.298354581028633
$ws.Range("C44").Value = 0.8146763841020004
$ws.Range("D44").Value = 0.3632612702442088
$ws.Range("E44").Value = 1.063894517653054
$ws.Range("F44").Value = 0.8709948986110732
$ws.Range("G44").Value = 0.2602534500261697
$ws.Range("H44").Value = 0.3656541826233735
$ws.Range("B45").Value = 0.4656763841019966
$ws.Range("C45").Value = 0.2632612702442003
$ws.Range("D45").Value = 0.8148945176530447
$ws.Range("E45").Value = 0.8219948986110523
$ws.Range("F45").Value = 0.3262534500261722
$ws.Range("G45").Value = 0.416654182623347
$ws.Range("H45").Value = 0.01016824066003608
$ws.Range("I45").Value = 0.2452723979283746
$ws.Range("B46").Value = -0.116738729755781
$ws.Range("C46").Value = 0.7148945176530219
$ws.Range("D46").Value = 0.6219948986110637
$ws.Range("E46").Value = 0.3112534500261574
$ws.Range("F46").Value = 0.5166541826233555
$ws.Range("G46").Value = 0.01016824066003608
$ws.Range("H46").Value = 0.3452723979283974
$ws.Range("B47").Value = 0.5148945176530333
$ws.Range("C47").Value = 0.6219948986110637
$ws.Range("D47").Value = 0.3112534500261574
$ws.Range("E47").Value = 0.5166541826233555
$ws.Range("F47").Value = 0.01016824066003608
$ws.Range("G47").Value = 0.3452723979283974
$ws.Range("B48").Value = 0.221994898611058
$ws.Range("C48").Value = 0.1602534500261754
$ws.Range("D48").Value = 0.3856541826233695
$ws.Range("E48").Value = 0.06116814066005816
$ws.Range("F48").Value = 0.505272397928394
$ws.Range("B49").Value = -0.1887465499738426
$ws.Range("C49").Value = 0.2166541826233441
$ws.Range("D49").Value = 0.01016824066003608
$ws.Range("E49").Value = 0.545272397928386
$ws.Range("B50").Value = 0.2366557702529377
$ws.Range("C50").Value = -0.1898317593399668
$ws.Range("D50").Value = 0.2052723979283826
$ws.Range("B51").Value = -0.3898340444052479
$ws.Range("C51").Value = 0.1852723979283866
$ws.Range("B52").Value = 0.1452758398526868

$ws.Range("J45").ClearContents()
$ws.Range("I46").ClearContents()
$ws.Range("H47").ClearContents()
$ws.Range("G48").ClearContents()
$ws.Range("F49").ClearContents()
$ws.Range("E50").ClearContents()
$ws.Range("D51").ClearContents()
$ws.Range("C52").ClearContents()
$ws.Range("B53").ClearContents()
